$d = $word.ActiveDocument

$replacements = @(
    @("35+13=48", "22-14=8"),
    @("92-49=43", "55-27=28"),
    @("67-22=45", "54-51=3"),
    @("73-67=6", "41-14=27"),
    @("49-34=15", "52-28=24"),
    @("20+14=34", "66+21=87"),
    @("25+72=97", "57-4=53"),
    @("26+8=34", "42+17=59"),
    @("26+4=30", "36-2=34"),
    @("61-47=14", "24+10=34"),
    @("14+76=90", "55-38=17"),
    @("7+42=49", "81-48=33"),
    @("26+17=43", "70-26=44"),
    @("20+26=46", "4+85=89"),
    @("41+51=92", "66+0=66"),
    @("62-40=22", "8+90=98"),
    @("1+98=99", "11+42=53"),
    @("44-29=15", "31+27=58"),
    @("11+88=99", "36-20=16"),
    @("27+18=45", "36-22=14"),
    @("55-24=31", "40+25=65"),
    @("88-26=62", "32-4=28"),
    @("32-19=13", "40-12=28"),
    @("14+59=73", "85-48=37"),
    @("4+15=19", "42-17=25"),
    @("64-15=49", "56-52=4"),
    @("53+23=76", "81-40=41"),
    @("29+50=79", "40-9=31"),
    @("93-32=61", "95-14=81"),
    @("72-21=51", "20+30=50"),
    @("25-9=16", "58+22=80"),
    @("75-64=11", "52+40=92"),
    @("13+39=52", "96-48=48"),
    @("58-35=23", "11-3=8"),
    @("90-21=69", "87-28=59"),
    @("69+18=87", "56-48=8"),
    @("98-27=71", "46-39=7"),
    @("89-61=28", "56+17=73"),
    @("82-53=29", "24-2=22"),
    @("31+43=74", "68-63=5"),
    @("58-30=28", "25-15=10"),
    @("91-31=60", "42+34=76"),
    @("68-36=32", "80+6=86"),
    @("29+36=65", "84-44=40"),
    @("2+33=35", "51-38=13"),
    @("62-29=33", "73-12=61"),
    @("50-22=28", "60-2=58"),
    @("14+29=43", "65+33=98"),
    @("23+10=33", "14-3=11"),
    @("64-31=33", "83-14=69"),
    @("13-7=6", "51-31=20"),
    @("10+33=43", "20+23=43"),
    @("17+69=86", "13+56=69"),
    @("56-23=33", "5+36=41"),
    @("68+0=68", "2+36=38"),
    @("35-10=25", "92-76=16"),
    @("46+7=53", "24-13=11"),
    @("12+43=55", "42-34=8"),
    @("34+35=69", "27+34=61"),
    @("38+43=81", "5+11=16"),
    @("66-23=43", "96-94=2"),
    @("40+12=52", "76-14=62"),
    @("32-10=22", "49-19=30"),
    @("98-69=29", "42+7=49"),
    @("9+2=11", "22+63=85"),
    @("65+28=93", "93-73=20"),
    @("0+47=47", "85-0=85"),
    @("83-74=9", "30+37=67"),
    @("18+60=78", "11+13=24"),
    @("27+48=75", "18-1=17"),
    @("23+7=30", "15+17=32"),
    @("78-21=57", "81-19=62"),
    @("52+36=88", "47+31=78"),
    @("4+66=70", "68+14=82"),
    @("87-53=34", "60-48=12"),
    @("25+70=95", "72-43=29"),
    @("36+10=46", "60-38=22"),
    @("65-49=16", "48-42=6"),
    @("4+84=88", "85-29=56"),
    @("34-27=7", "76+6=82"),
    @("9+78=87", "33+46=79"),
    @("82-80=2", "15+77=92"),
    @("95-60=35", "73-54=19"),
    @("59-33=26", "93-0=93"),
    @("12+38=50", "79+13=92"),
    @("22+26=48", "18+62=80"),
    @("99-18=81", "97-22=75"),
    @("45+35=80", "40-14=26"),
    @("1+70=71", "70-10=60"),
    @("32+41=73", "74-37=37"),
    @("2+40=42", "33+66=99"),
    @("75+18=93", "83+16=99"),
    @("80-56=24", "85-83=2"),
    @("88-30=58", "54+36=90"),
    @("43-37=6", "65+15=80"),
    @("78+7=85", "19+77=96"),
    @("34-34=0", "1+62=63"),
    @("80+11=91", "40-35=5"),
    @("98-3=95", "95-73=22"),
    @("68-30=38", "72-54=18")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2)
}
